$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 113 (so they become the new rows 114 and 115),
# pushing the existing rows 114-135 down to 116-137.
$ws.Rows.Item(114).Resize(2).Insert()

# Row 114 (new)
$ws.Cells.Item(114, 1).Value = 1
$ws.Cells.Item(114, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(114, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(114, 4).Value = 44855
$ws.Cells.Item(114, 5).Value = 15
$ws.Cells.Item(114, 6).Value = "Fruta"
$ws.Cells.Item(114, 7).Value = 100106
$ws.Cells.Item(114, 8).Value = "Oleaginosos"
$ws.Cells.Item(114, 9).Value = 100106002
$ws.Cells.Item(114, 10).Value = "Palta"
$ws.Cells.Item(114, 11).Value = "Fuerte"
$ws.Cells.Item(114, 12).Value = "Segunda"
$ws.Cells.Item(114, 13).Value = 600
$ws.Cells.Item(114, 14).Value = 9000
$ws.Cells.Item(114, 15).Value = 10000
$ws.Cells.Item(114, 16).Value = 9583
$ws.Cells.Item(114, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(114, 18).Value = "Perú"
$ws.Cells.Item(114, 19).Value = 958
$ws.Cells.Item(114, 20).Value = 10

# Row 115 (new)
$ws.Cells.Item(115, 1).Value = 1
$ws.Cells.Item(115, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(115, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(115, 4).Value = 44855
$ws.Cells.Item(115, 5).Value = 15
$ws.Cells.Item(115, 6).Value = "Fruta"
$ws.Cells.Item(115, 7).Value = 100106
$ws.Cells.Item(115, 8).Value = "Oleaginosos"
$ws.Cells.Item(115, 9).Value = 100106002
$ws.Cells.Item(115, 10).Value = "Palta"
$ws.Cells.Item(115, 11).Value = "Hass"
$ws.Cells.Item(115, 12).Value = "Tercera"
$ws.Cells.Item(115, 13).Value = 900
$ws.Cells.Item(115, 14).Value = 23000
$ws.Cells.Item(115, 15).Value = 24000
$ws.Cells.Item(115, 16).Value = 23556
$ws.Cells.Item(115, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(115, 18).Value = "Perú"
$ws.Cells.Item(115, 19).Value = 2356
$ws.Cells.Item(115, 20).Value = 10

# Apply the same date format/style as column D in the other rows to the new D cells
$ws.Range("D114:D115").NumberFormat = $ws.Range("D113").NumberFormat
